$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct typo: "6.5." -> "7.4." in the args_pre description (C2)
$ws.Range("C2").Value = "Path to a Yaml file or dictionary containing preprocessing configuration information.`nRefer to 7.4. for details."

# Correct typo: add missing period after ".txt" in the save_name description (C5)
$ws.Range("C5").Value = "Name for resulting files.`nNumpy files will be saved under {save_dir}/{save_name}_npy directory.`nText file will be saved in {save_dir}/{save_name}.txt.`nIf it is not provided, it is set to the basename of data_dir."

# Move the active selection from C3 to C2
[void]$ws.Range("C2").Select()
